$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value without altering its cell style,
# even when the text looks like a number (e.g. "0.9999", "1.0000").
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "29.481.12"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.880.08"
$ws.Range("E3").Value = "  +1.28%  "
Set-TextValue $ws.Range("D4") "0.9999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +1.71%  "
Set-TextValue $ws.Range("D6") "240.44"
$ws.Range("E6").Value = "  +0.94%  "
Set-TextValue $ws.Range("D7") "1.0000"
$ws.Range("E7").Value = "  -0.03%  "
Set-TextValue $ws.Range("D8") "0.07889"
$ws.Range("E8").Value = "  -1.79%  "
Set-TextValue $ws.Range("D9") "0.3123"
$ws.Range("E9").Value = "  +3.36%  "
Set-TextValue $ws.Range("D10") "25.22"
$ws.Range("E10").Value = "  +7.46%  "
Set-TextValue $ws.Range("D11") "0.08248"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "1.880.06"
$ws.Range("E12").Value = "  +1.38%  "
Set-TextValue $ws.Range("D13") "0.7309"
$ws.Range("E13").Value = "  +3.99%  "
Set-TextValue $ws.Range("D14") "5.300"
$ws.Range("E14").Value = "  +2.16%  "
Set-TextValue $ws.Range("D15") "91.31"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "29.622.20"
$ws.Range("E16").Value = "  +1.72%  "
Set-TextValue $ws.Range("D17") "5.949"
$ws.Range("E17").Value = "  +2.80%  "
Set-TextValue $ws.Range("D18") "248.93"
$ws.Range("E18").Value = "  +4.98%  "
Set-TextValue $ws.Range("D19") "0.000007899"
$ws.Range("E19").Value = "  -0.92%  "
Set-TextValue $ws.Range("D20") "13.32"
$ws.Range("E20").Value = "  +0.83%  "
Set-TextValue $ws.Range("D21") "0.9996"
$ws.Range("E21").Value = "  +0.00%  "
Set-TextValue $ws.Range("D22") "8.006"
$ws.Range("E22").Value = "  +7.25%  "
Set-TextValue $ws.Range("D23") "1.001"
$ws.Range("E23").Value = "  +0.05%  "
Set-TextValue $ws.Range("D24") "0.1566"
$ws.Range("E24").Value = "  +9.67%  "
Set-TextValue $ws.Range("D25") "164.03"
$ws.Range("E25").Value = "  +0.77%  "
Set-TextValue $ws.Range("D26") "9.062"
$ws.Range("E26").Value = "  +1.79%  "
Set-TextValue $ws.Range("D27") "18.37"
$ws.Range("E27").Value = "  +1.64%  "
Set-TextValue $ws.Range("D28") "1.366"
$ws.Range("E28").Value = "  -3.57%  "
Set-TextValue $ws.Range("D29") "1.490"
$ws.Range("E29").Value = "  +0.89%  "
Set-TextValue $ws.Range("D30") "4.392"
$ws.Range("E30").Value = "  +0.85%  "
Set-TextValue $ws.Range("D31") "4.154"
$ws.Range("E31").Value = "  +3.24%  "
Set-TextValue $ws.Range("D32") "0.05294"
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  +1.19%  "
Set-TextValue $ws.Range("D34") "1.207"
$ws.Range("E34").Value = "  +4.39%  "
Set-TextValue $ws.Range("D35") "0.7251"
$ws.Range("E35").Value = "  +1.54%  "
Set-TextValue $ws.Range("D36") "2.676"
$ws.Range("E36").Value = "  +1.18%  "
Set-TextValue $ws.Range("D37") "0.01868"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").Value = "1.271.94"
$ws.Range("E38").Value = "  +12.84%  "
$ws.Range("E39").Value = "  +0.19%  "
Set-TextValue $ws.Range("D40") "0.9063"
$ws.Range("E40").Value = "  -3.28%  "
Set-TextValue $ws.Range("D41") "73.82"
$ws.Range("E41").Value = "  +5.39%  "
Set-TextValue $ws.Range("D42") "6.117"
$ws.Range("E42").Value = "  +3.23%  "
Set-TextValue $ws.Range("D43") "103.96"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("E44").Value = "  -0.01%  "
Set-TextValue $ws.Range("D45") "0.5332"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "1.772"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D47") "2.925"
$ws.Range("E47").Value = "  +13.20%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.00000000120"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D49") "0.4339"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "9.292"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D51") "7.092"
$ws.Range("E51").Value = "  +2.16%  "
